$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.926.17"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "1.644.53"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  -0.69%  "

$ws.Range("D5").Value = "216.99"

$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("E8").Value = "  +1.15%  "

$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").Value = "19.85"
$ws.Range("E10").Value = "  +4.02%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "1.873.90"
$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("D13").Value = "1.654.76"
$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("E14").Value = "  +0.23%  "

$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").Value = "66.27"
$ws.Range("E16").Value = "  +2.79%  "

$ws.Range("D17").Value = "26.945.86"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").Value = "219.00"
$ws.Range("E19").Value = "  +3.62%  "

$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").Value = "4.39"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  +7.30%  "

$ws.Range("E23").Value = "  +6.08%  "

$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("D25").Value = "145.72"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("D27").Value = "7.39"
$ws.Range("E27").Value = "  +4.45%  "

$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("D29").Value = "15.86"
$ws.Range("E29").Value = "  +1.86%  "

$ws.Range("D30").Value = "0.0512"
$ws.Range("E30").Value = "  +1.84%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("D34").Value = "1.56"
$ws.Range("E34").Value = "  +2.09%  "

$ws.Range("D35").Value = "2.44"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D36").Value = "1.247.30"
$ws.Range("E36").Value = "  -1.71%  "

$ws.Range("D37").Value = "0.0175"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").Value = "0.538"
$ws.Range("E38").Value = "  +1.92%  "

$ws.Range("E39").Value = "  +3.91%  "

$ws.Range("E40").Value = "  -0.71%  "

$ws.Range("E41").Value = "  +0.81%  "

$ws.Range("D42").Value = "5.35"
$ws.Range("E42").Value = "  +1.56%  "

$ws.Range("D43").Value = "1.786.44"
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").Value = "2.11"
$ws.Range("E44").Value = "  -3.52%  "

$ws.Range("D45").Value = "61.00"
$ws.Range("E45").Value = "  +1.45%  "

$ws.Range("D46").Value = "91.57"
$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("E47").Value = "  +1.33%  "

$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  -0.86%  "

$ws.Range("D50").Value = "0.0973"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("D51").Value = "7.56"
$ws.Range("E51").Value = "  +0.43%  "
